# Daily attendance processing - reorder "Recorded By" (column G) entries.
# For every data row, if the comma-separated list of recorders contains an
# entry equal to "System" (case-insensitive) and has more than one entry,
# reverse the order of the entries in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "

    if ($parts.Count -gt 1) {
        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Trim().ToLower() -eq "system") {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newText = [string]::Join(", ", $reversed)
            $cell.Value2 = $newText
        }
    }
}
